# Refresh the cryptos price/volume table (and two row re-ranks + one
# row replacement) to match the latest GitHub Actions scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price column holds values like "68.046.24" or "1.00" that Excel
# would otherwise auto-coerce to numbers/dates; force it to Text first
# so COM stores the exact literal strings the source data uses.
$ws.Range("D2:D51").NumberFormat = "@"

# Row 2 - Bitcoin
$ws.Range("D2").Value = '68.046.24'
$ws.Range("E2").Value = '  -0.43%  '

# Row 3 - Ethereum
$ws.Range("D3").Value = '3.617.56'
$ws.Range("E3").Value = '  -1.60%  '

# Row 4 - TetherUSD
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  +0.15%  '

# Row 5 - BNB
$ws.Range("D5").Value = '586.99'
$ws.Range("E5").Value = '  -2.22%  '

# Row 6 - Solana
$ws.Range("D6").Value = '193.92'
$ws.Range("E6").Value = '  +0.33%  '

# Row 7 - LidoStakedEther
$ws.Range("D7").Value = '3.613.01'
$ws.Range("E7").Value = '  -1.51%  '

# Row 8 - XRP
$ws.Range("D8").Value = '0.620'
$ws.Range("E8").Value = '  -0.65%  '

# Row 9 - USDC
$ws.Range("E9").Value = '  +0.22%  '

# Row 10 - Cardano
$ws.Range("D10").Value = '0.682'
$ws.Range("E10").Value = '  -2.95%  '

# Row 11 - Dogecoin
$ws.Range("D11").Value = '0.152'
$ws.Range("E11").Value = '  -1.36%  '

# Row 12 - Avalanche
$ws.Range("D12").Value = '55.79'
$ws.Range("E12").Value = '  -3.33%  '

# Row 13 - ShibaInu
$ws.Range("D13").Value = '0.0000293'
$ws.Range("E13").Value = '  +6.92%  '

# Row 14 - Polkadot
$ws.Range("D14").Value = '10.02'
$ws.Range("E14").Value = '  -2.46%  '

# Row 15 - WrappedliquidstakedEther2.0
$ws.Range("D15").Value = '4.194.75'
$ws.Range("E15").Value = '  -1.33%  '

# Row 16 - WrappedEther
$ws.Range("D16").Value = '3.614.41'
$ws.Range("E16").Value = '  -1.55%  '

# Row 17 - TRON
$ws.Range("E17").Value = '  -0.59%  '

# Row 18 - Uniswap
$ws.Range("D18").Value = '12.56'
$ws.Range("E18").Value = '  -0.99%  '

# Row 19 - WrappedBTC
$ws.Range("B19").Value = 'WrappedBTC'
$ws.Range("C19").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D19").Value = '67.916.95'
$ws.Range("E19").Value = '  -0.17%  '

# Row 20 - Chainlink
$ws.Range("B20").Value = 'Chainlink'
$ws.Range("C20").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D20").Value = '18.53'
$ws.Range("E20").Value = '  -2.34%  '

# Row 21 - Polygon
$ws.Range("E21").Value = '  -2.92%  '

# Row 22 - BitcoinCash
$ws.Range("D22").Value = '404.82'
$ws.Range("E22").Value = '  -0.79%  '

# Row 23 - RenderToken
$ws.Range("D23").Value = '13.52'
$ws.Range("E23").Value = '  +22.70%  '

# Row 24 - PancakeSwap
$ws.Range("D24").Value = '4.27'
$ws.Range("E24").Value = '  -3.98%  '

# Row 25 - Litecoin
$ws.Range("D25").Value = '86.08'
$ws.Range("E25").Value = '  -2.80%  '

# Row 26 - ImmutableX
$ws.Range("D26").Value = '2.95'
$ws.Range("E26").Value = '  -0.41%  '

# Row 27 - InternetComputer(DFINITY)
$ws.Range("D27").Value = '12.67'
$ws.Range("E27").Value = '  -0.32%  '

# Row 28 - Toncoin
$ws.Range("D28").Value = '3.95'
$ws.Range("E28").Value = '  +6.17%  '

# Row 29 - LEO
$ws.Range("D29").Value = '6.12'
$ws.Range("E29").Value = '  +0.80%  '

# Row 30 - NEARProtocol
$ws.Range("D30").Value = '8.32'
$ws.Range("E30").Value = '  +15.30%  '

# Row 31 - Filecoin
$ws.Range("D31").Value = '9.22'
$ws.Range("E31").Value = '  -2.01%  '

# Row 32 - EthereumClassic
$ws.Range("D32").Value = '31.53'
$ws.Range("E32").Value = '  -1.92%  '

# Row 33 - Bittensor
$ws.Range("D33").Value = '684.98'
$ws.Range("E33").Value = '  +11.60%  '

# Row 34 - Cosmos
$ws.Range("D34").Value = '12.27'
$ws.Range("E34").Value = '  -0.82%  '

# Row 35 - Hedera
$ws.Range("E35").Value = '  +1.20%  '

# Row 36 - OKB
$ws.Range("D36").Value = '64.58'
$ws.Range("E36").Value = '  -3.35%  '

# Row 37 - InjectiveProtocol
$ws.Range("D37").Value = '42.60'
$ws.Range("E37").Value = '  -4.30%  '

# Row 38 - TheGraph
$ws.Range("D38").Value = '0.422'
$ws.Range("E38").Value = '  +6.33%  '

# Row 39 - Dai
$ws.Range("B39").Value = 'Dai'
$ws.Range("C39").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D39").Value = '1.00'
$ws.Range("E39").Value = '  +0.23%  '

# Row 40 - PEPE
$ws.Range("B40").Value = 'PEPE'
$ws.Range("C40").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D40").Value = '0.0₃0788'
$ws.Range("E40").Value = '  +0.53%  '

# Row 41 - Fetch.AI
$ws.Range("D41").Value = '2.99'
$ws.Range("E41").Value = '  +17.73%  '

# Row 42 - ThetaToken
$ws.Range("D42").Value = '3.14'
$ws.Range("E42").Value = '  +7.53%  '

# Row 43 - Maker
$ws.Range("D43").Value = '3.191.41'
$ws.Range("E43").Value = '  +13.96%  '

# Row 44 - Kaspa
$ws.Range("D44").Value = '0.134'
$ws.Range("E44").Value = '  -1.52%  '

# Row 45 - FirstDigitalUSD
$ws.Range("D45").Value = '0.998'
$ws.Range("E45").Value = '  -0.09%  '

# Row 46 - VeChain
$ws.Range("D46").Value = '0.0421'
$ws.Range("E46").Value = '  -1.61%  '

# Row 47 - Stellar
$ws.Range("D47").Value = '0.132'
$ws.Range("E47").Value = '  -3.11%  '

# Row 48 - THORChain
$ws.Range("D48").Value = '8.82'
$ws.Range("E48").Value = '  -1.50%  '

# Row 49 - Monero
$ws.Range("D49").Value = '143.74'
$ws.Range("E49").Value = '  +0.29%  '

# Row 50 - ApeXProtocol
$ws.Range("D50").Value = '3.08'
$ws.Range("E50").Value = '  -3.32%  '

# Row 51 - Stacks
$ws.Range("B51").Value = 'Stacks'
$ws.Range("C51").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D51").Value = '2.78'
$ws.Range("E51").Value = '  +2.42%  '

# Restore the default cell style so the text-format override above
# does not leave a stray NumberFormat on these cells.
$ws.Range("D2:D51").Style = "Normal"
